$d = $word.ActiveDocument

# --- Change 1: fill in Angela's "Work completed" cell (table 1, row 3, col 2) ---
$t = $d.Tables.Item(1)
$cell = $t.Cell(3, 2)
$newText = "Continue to work on the website, researched and tried multiple codes for adding working calendar for booking system. Participated in creating of presentation."

# Insert the text plus a sentinel character so the bookmark can be placed by
# splitting a run (collapsed ranges that fall exactly on a paragraph/run
# boundary don't anchor correctly), then trim the sentinel back out.
$r = $cell.Range
$r.InsertBefore($newText + "X")

$cell2 = $t.Cell(3, 2)
$rng2 = $cell2.Range
$bmPos = $rng2.Start + $newText.Length
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

$delRng = $d.Range($bmPos, $bmPos + 1)
$delRng.Delete()

# --- Change 2: merge the "YouT" / "ube" runs (split by the old _GoBack
#     bookmark) back into a single contiguous run of text ---
$null = $d.Content.Find.Execute("Watched a lot of video on YouT*ube to help", $true, $false, $true, $false, $false, $true, 1, $false, "Watched a lot of video on YouTube to help", 2)
